# 11th botany TM U02 Q&A added
# The title bar (Rectangle 1) on every question/answer slide grows taller
# (to fit its two-line heading) and the body textbox (TextBox 2) slides
# down so it continues to sit directly below the title bar.
#
# EMU -> point conversion is exact (1 pt = 12700 EMU), but the COM layer's
# internal float rounds a couple of these conversions down by a single
# EMU, so a hair (0.00001 pt ~= 0.127 EMU) is added on top of the exact
# point value to make sure it rounds back to the intended EMU value.

$p = $ppt.ActivePresentation

for ($i = 2; $i -le 43; $i++) {
    $s = $p.Slides.Item($i)

    $rect = $s.Shapes.Item(1)
    $rect.Top = 21.60001
    $rect.Height = 57.60001

    $textBox = $s.Shapes.Item(2)
    $textBox.Top = 79.20001
}
